$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each rewritten cell to a Text format *before* assigning its new
# string so Excel keeps the exact literal (e.g. "99.10", "0.0840") instead
# of silently reinterpreting it as a number and dropping trailing zeros.
# (Applied per-cell, not as a bulk/union range, since multi-area Range
# property assignment only reliably touches the first area here.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.108.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.57%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.26"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.38"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.50%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.585.61"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.249.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.962.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.32%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.40%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.07"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.91"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0840"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.82%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.38"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.53"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -10.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.98"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -11.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0309"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.707.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.19"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.78"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.03"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.50%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.24%  "
